$d = $word.ActiveDocument
$d.Content.Find.Execute("custom on-demand builds of Material-UI", $true, $false, $false, $false, $false,
                         $true, 1, $false, "custom on-demand build of Material-UI", 2)
